$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.636.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5335"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3176"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06962"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7720"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07826"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.867.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.052"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007990"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.671.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.088.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.654"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.032"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.387"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.215"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.705"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.330"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08768"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04867"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7405"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.888"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.106"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.360"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01747"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4837"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9085"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.715"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4212"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1252"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05827"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.08%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "

